$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 held per-column template placeholders that pulled values from a
# linked child sheet (e.g. ='{child:sheetname}'!H5). Those placeholders are
# being retired in favor of the same "{empty}" placeholder/format already
# used by D7 (and H7/L7's neighbors), so copy D7's format onto them and set
# their text to match.
$src = $ws.Range("D7")
$emptyValue = $src.Value2

$targets = @("E7", "F7", "G7", "I7", "J7", "K7", "M7", "N7")

$src.Copy() | Out-Null
foreach ($addr in $targets) {
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = 0

foreach ($addr in $targets) {
    $ws.Range($addr).Value = $emptyValue
}

# Move the active selection to N7, matching where the author last edited.
$ws.Range("N7").Select() | Out-Null
